$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SP25092022103921", 100, "PANADOL STRIP 10", 2, 60),
    @("SP25092022103921", 103, "NUROFEN STRIP 25", 4, 400),
    @("SP25092022104544", 100, "PANADOL STRIP 10", 2, 60),
    @("SP25092022104544", 103, "NUROFEN STRIP 25", 2, 100),
    @("SP25092022105145", 100, "PANADOL STRIP 10", 5, 375),
    @("SP25092022105145", 103, "NUROFEN STRIP 25", 5, 625),
    @("SP25092022105547", 100, "PANADOL STRIP 10", 2, 60),
    @("SP25092022105547", 103, "NUROFEN STRIP 25", 4, 400),
    @("SP25092022110035", 100, "PANADOL STRIP 10", 2, 60),
    @("SP25092022110035", 103, "NUROFEN STRIP 25", 2, 100),
    @("SP25092022110524", 100, "PANADOL STRIP 10", 2, 60),
    @("SP25092022110524", 102, "NUROFEN STRIP 15", 2, 48),
    @("SP25092022110937", 100, "PANADOL STRIP 10", 2, 60),
    @("SP25092022110937", 102, "NUROFEN STRIP 15", 1, 12)
)

$startRow = 33
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$wb.Save()
